$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the "vbat" description cell (E4 was rich text with 3 runs) to a
# single plain-text run with the same concatenated text.
$ws.Range("E4").Value = "12-bit ADC reading of battery measurement N[11:0], reference R=0: 1.5 V, R=1: 2.5 V"

# Add new rows describing the BMP180 temperature and pressure sensor fields.
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "s"
$ws.Range("D8").Value = "bmp180_temperature"
$ws.Range("E8").Value = "BMP180 temperature reading in 0.1 degree Celsius"

$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "s"
$ws.Range("D9").Value = "bmp180_pressure"
$ws.Range("E9").Value = "BMP180 pressure reading in Pascal"

# Update the default column width / active selection to match the state the
# workbook was saved in after this edit.
$ws.Range("B9").Select()

$wb.Save()
